# Update "想去人数" (interest count) values on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - rows 2-12
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 755
$ws1.Range("F3").Value = 12
$ws1.Range("F4").Value = 49
$ws1.Range("F5").Value = 19
$ws1.Range("F6").Value = 256
$ws1.Range("F7").Value = 3388
$ws1.Range("F9").Value = 4064
$ws1.Range("F11").Value = 1019
$ws1.Range("F12").Value = 41

# Sheet "全部类型" (all types) - rows 2-13 (offset by one extra row vs 展览)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 755
$ws4.Range("F3").Value = 12
$ws4.Range("F4").Value = 49
$ws4.Range("F5").Value = 19
$ws4.Range("F7").Value = 256
$ws4.Range("F8").Value = 3388
$ws4.Range("F10").Value = 4064
$ws4.Range("F12").Value = 1019
$ws4.Range("F13").Value = 41
